$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the fill-down of the B column formula (A{row}-(((A{row}+45)*(A{row}+45))*$D$71))
# from B118 down through B131, matching column A which already goes to row 131.
$ws.Range("B119:B131").Formula = '=A119-(((A119+45)*(A119+45))*$D$71)'

# Match the existing centered-number style used by the rest of column B/A (style index 1).
$ws.Range("B119:B131").HorizontalAlignment = -4108

# Update the view: scrolled down so rows ~114+ are visible, with B118:B131 selected.
$ws.Range("B118:B131").Select()
